$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2 through 250 change from serial date 45203 (2023-10-04)
# to serial date 45205 (2023-10-06).
for ($r = 2; $r -le 250; $r++) {
    $ws.Cells.Item($r, 3).Value = 45205
}
